# PaymentsELHardCoded.xlsx -- "Fixed RAD Test Cases and Data. Added VLink
# Smoke Test Cases and Data." re-run of the RAD recorder against Sheet1:
#   - refreshed the "Date" column (B2:B14) timestamps with a new run's log
#   - left the cursor/selection on a different cell (P5) after the run
#   - a sensitivity-label add-in stamped a footer onto the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B ("Date") timestamps refreshed for the new test run --------
$ws.Range("B2").Value  = "Wed Jan 04 19:01:22 EST 2023"
$ws.Range("B3").Value  = "Wed Jan 04 19:02:11 EST 2023"
$ws.Range("B4").Value  = "Wed Jan 04 19:02:59 EST 2023"
$ws.Range("B5").Value  = "Wed Jan 04 19:03:48 EST 2023"
$ws.Range("B6").Value  = "Wed Jan 04 19:05:15 EST 2023"
$ws.Range("B7").Value  = "Wed Jan 04 19:06:07 EST 2023"
$ws.Range("B8").Value  = "Wed Jan 04 19:06:56 EST 2023"
$ws.Range("B9").Value  = "Wed Jan 04 19:07:44 EST 2023"
$ws.Range("B10").Value = "Wed Jan 04 19:08:32 EST 2023"
$ws.Range("B11").Value = "Wed Jan 04 19:09:23 EST 2023"
$ws.Range("B12").Value = "Wed Jan 04 19:10:11 EST 2023"
$ws.Range("B13").Value = "Wed Jan 04 19:11:01 EST 2023"
$ws.Range("B14").Value = "Wed Jan 04 19:11:49 EST 2023"

# --- Selection left on P5 when the workbook was saved -------------------
$ws.Range("P5").Select()

# --- Footer stamped by the org's sensitivity-label add-in ---------------
# ("&C" center-section code, a literal CR, then font/size/color codes
# around the word "Public".)
$ps = $ws.PageSetup
$ps.CenterFooter = "`r&1#&`"Calibri`"&10&K000000 Public "
